$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.280.69'
$ws.Range('E2').Value = '  +0.38%  '

$ws.Range('D3').Value = '3.492.41'
$ws.Range('E3').Value = '  +1.05%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '588.14'
$ws.Range('E5').Value = '  -0.65%  '

$ws.Range('D6').Value = '177.72'
$ws.Range('E6').Value = '  -0.62%  '

$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  -1.50%  '

$ws.Range('D9').Value = '3.490.00'
$ws.Range('E9').Value = '  +1.20%  '

$ws.Range('D10').Value = '0.134'
$ws.Range('E10').Value = '  -2.00%  '

$ws.Range('D11').Value = '6.92'
$ws.Range('E11').Value = '  -0.25%  '

$ws.Range('D12').Value = '0.424'
$ws.Range('E12').Value = '  -1.55%  '

$ws.Range('D13').Value = '4.080.51'
$ws.Range('E13').Value = '  +0.89%  '

$ws.Range('D14').Value = '30.73'
$ws.Range('E14').Value = '  -2.63%  '

$ws.Range('E15').Value = '  -0.91%  '

$ws.Range('D16').Value = '67.214.72'
$ws.Range('E16').Value = '  +0.45%  '

$ws.Range('E17').Value = '  -0.93%  '

$ws.Range('D18').Value = '3.473.44'
$ws.Range('E18').Value = '  +0.56%  '

$ws.Range('D19').Value = '6.06'
$ws.Range('E19').Value = '  -2.73%  '

$ws.Range('D20').Value = '13.96'
$ws.Range('E20').Value = '  -1.03%  '

$ws.Range('D21').Value = '384.90'
$ws.Range('E21').Value = '  -0.71%  '

$ws.Range('D22').Value = '7.90'
$ws.Range('E22').Value = '  -0.13%  '

$ws.Range('E23').Value = '  +2.35%  '

$ws.Range('D24').Value = '73.09'
$ws.Range('E24').Value = '  +1.21%  '

$ws.Range('D25').Value = '5.81'
$ws.Range('E25').Value = '  +0.93%  '

$ws.Range('E26').Value = '  +0.33%  '

$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').Value = '  +0.84%  '

$ws.Range('D28').Value = '9.94'
$ws.Range('E28').Value = '  -3.10%  '

$ws.Range('E29').Value = '  +1.74%  '

$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.27%  '

$ws.Range('D31').Value = '24.42'
$ws.Range('E31').Value = '  +4.75%  '

$ws.Range('D32').Value = '5.94'
$ws.Range('E32').Value = '  -3.36%  '

$ws.Range('E33').Value = '  -1.85%  '

$ws.Range('E34').Value = '  -3.79%  '

$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('D36').Value = '7.22'
$ws.Range('E36').Value = '  -1.23%  '

$ws.Range('E37').Value = '  +0.43%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '30.01'
$ws.Range('E38').Value = '  +15.59%  '

$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '160.96'
$ws.Range('E39').Value = '  -0.38%  '

$ws.Range('D40').Value = '0.898'
$ws.Range('E40').Value = '  +2.82%  '

$ws.Range('D41').Value = '1.81'
$ws.Range('E41').Value = '  -2.71%  '

$ws.Range('D42').Value = '2.63'
$ws.Range('E42').Value = '  -5.29%  '

$ws.Range('D43').Value = '4.54'
$ws.Range('E43').Value = '  -2.35%  '

$ws.Range('D44').Value = '6.51'
$ws.Range('E44').Value = '  -4.30%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.711.80'
$ws.Range('E45').Value = '  -2.11%  '

$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0702'
$ws.Range('E46').Value = '  -2.20%  '

$ws.Range('D47').Value = '40.78'
$ws.Range('E47').Value = '  -0.70%  '

$ws.Range('D48').Value = '24.81'
$ws.Range('E48').Value = '  -4.41%  '

$ws.Range('D49').Value = '0.0296'
$ws.Range('E49').Value = '  -0.13%  '

$ws.Range('D50').Value = '319.07'
$ws.Range('E50').Value = '  -2.07%  '

$ws.Range('E51').Value = '  -1.31%  '
